# Bugfix in steiner tree rebuilding alternative random node including tests
# (commented out for now) run_opt script.
#
# This adds a new results sheet "GRASP (randnode2_rec_seed)" (copied from the
# "GRASP (randnode_rec_seed)" sheet so it inherits the same layout, styles,
# formulas and conditional formatting), fills it with the new run's data,
# and fixes up a value on the "GRASP (randnode_rec_seed)" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Fix bugged value on "GRASP (randnode_rec_seed)" (4th sheet) -------
$wsRandnode = $wb.Worksheets.Item(4)
$wsRandnode.Range("H11").Value = 109672

# Sheet is no longer the active/selected one; selection moves to D8.
$wsRandnode.Activate()
$wsRandnode.Range("D8").Select()

# --- 2. Add the new sheet as a copy of "GRASP (randnode_rec_seed)" --------
$wsRandnode.Copy($null, $wsRandnode)
$newWs = $wb.Worksheets.Item(5)
$newWs.Name = "GRASP (randnode2_rec_seed)"

# --- 3. Populate the new sheet with the new run's data (B2:K11) -----------
$values = @(
    @(20,42783759,132894610,7616,39144,22068,102794,95859,34316,80775),
    @(20,42783759,132987802,7375,39254,22068,105772,96026,34316,80775),
    @(20,42630715,132894610,7375,38561,22068,106214,95781,34316,80775),
    @(20,42085075,132894610,7407,38826,22068,107208,97304,34316,80775),
    @(20,42630715,132330304,7382,38485,22068,108199,102892,34316,80236),
    @(20,42630715,132245121,7610,38334,22068,108760,102005,34316,80775),
    @(20,42630715,132245121,7610,38334,22068,108760,102005,34316,80775),
    @(20,42447877,131294566,7307,40154,22068,106314,102892,34316,80775),
    @(20,42630715,131294566,7503,40234,22068,108986,102892,34316,80775),
    @(20,42783759,131022310,7434,42204,22068,107522,102892,34316,80775)
)

$arr = New-Object 'object[,]' 10,10
for ($r = 0; $r -lt 10; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt 10; $c++) {
        $arr[$r,$c] = $row[$c]
    }
}
$newWs.Range("B2:K11").Value = $arr

# --- 4. New sheet becomes the active tab / selected sheet ------------------
$newWs.Activate()
$newWs.Range("I9").Select()
